# "add more res in บอทโง่ tag"
# Append one more response row for the existing "บอทโง่" tag, mirroring the
# tag value already stored in A23 and adding a fresh response string in B24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tagValue = $ws.Range("A23").Value2
$ws.Range("A24").Value2 = $tagValue
$ws.Range("B24").Value2 = "ซักหมัดป้ะ"

# Match the author's final selection/viewport in the saved file.
[void]$ws.Range("B28").Select()
